$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the old DARWIN REALES CASTRO worker row (old row 17) -- shifts rows
# below it up by one and drops the now-unused shared strings for that row.
$ws.Rows("17").Delete()

# Update the aggregate "VALOR MORA" total.
$ws.Range("E11").Value = 1547

# Update worker / period counts.
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
